$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" labels (E16:E22) so that periods run in ascending
# order from 2203 (oldest) at the top to 2209 (most recent) at the bottom.
$ws.Range("E16").Value = "2203"
$ws.Range("E17").Value = "2204"
$ws.Range("E18").Value = "2205"
$ws.Range("E19").Value = "2206"
$ws.Range("E20").Value = "2207"
$ws.Range("E21").Value = "2208"
$ws.Range("E22").Value = "2209"

# "Valor Mora" (F column): the partial-period value (34666) now belongs to
# period 2209 (row 22); every other period uses the full value 40000.
$ws.Range("F16").Value = 40000
$ws.Range("F17").Value = 40000
$ws.Range("F18").Value = 40000
$ws.Range("F19").Value = 40000
$ws.Range("F20").Value = 40000
$ws.Range("F21").Value = 40000
$ws.Range("F22").Value = 34666

# "Salario Basico" (G column): updated from 1508000 to 1000000 for every row.
$ws.Range("G16:G22").Value = 1000000
